$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a 4th row of formulas mirroring the existing pattern in rows 1-3:
# col A holds a standalone formula, cols B:C share one formula (B is the
# "master" with the shared ref/definition, C just points at it).
$ws.Range("A4").Formula = "=E4*F4"
$ws.Range("B4:C4").Formula = "=F4*G4"
